# Natmi following Dr Hou advice
# Update row 2 (ECs -> ECs target) values, rewrite row 3 (ECs -> FAPs target),
# and add two new rows (ECs -> M1, ECs -> sCs) to the LR-pairs table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2: Target cluster changes from "FAPs" to "ECs" plus refreshed metrics ----
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.1847786666666667
$ws.Cells.Item(2, 8).Value = 0.5543360000000001
$ws.Cells.Item(2, 9).Value = 1
$ws.Cells.Item(2, 10).Value = 1
$ws.Cells.Item(2, 11).Value = 1
$ws.Cells.Item(2, 12).Value = 0.3333333333333333
$ws.Cells.Item(2, 13).Value = 0.02608566666666666
$ws.Cells.Item(2, 14).Value = 0.07825699999999999
$ws.Cells.Item(2, 15).Value = 0.007055522672798636
$ws.Cells.Item(2, 16).Value = 0.007055522672798637
$ws.Cells.Item(2, 17).Value = 0.004820074705777777
$ws.Cells.Item(2, 18).Value = 0.043380672352
$ws.Cells.Item(2, 19).Value = 0.007055522672798636
$ws.Cells.Item(2, 20).Value = 0.007055522672798637

# ---- Row 3: Target cluster stays "FAPs" but metrics refreshed ----
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vip"
$ws.Cells.Item(3, 3).Value = "Vipr2"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.1847786666666667
$ws.Cells.Item(3, 8).Value = 0.5543360000000001
$ws.Cells.Item(3, 9).Value = 1
$ws.Cells.Item(3, 10).Value = 1
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.6658376666666667
$ws.Cells.Item(3, 14).Value = 1.997513
$ws.Cells.Item(3, 15).Value = 0.1800924934601381
$ws.Cells.Item(3, 16).Value = 0.1800924934601381
$ws.Cells.Item(3, 17).Value = 0.1230325962631111
$ws.Cells.Item(3, 18).Value = 1.107293366368
$ws.Cells.Item(3, 19).Value = 0.1800924934601381
$ws.Cells.Item(3, 20).Value = 0.1800924934601381

# ---- Row 4 (new): Target cluster "M1" ----
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vip"
$ws.Cells.Item(4, 3).Value = "Vipr2"
$ws.Cells.Item(4, 4).Value = "M1"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.1847786666666667
$ws.Cells.Item(4, 8).Value = 0.5543360000000001
$ws.Cells.Item(4, 9).Value = 1
$ws.Cells.Item(4, 10).Value = 1
$ws.Cells.Item(4, 11).Value = 2
$ws.Cells.Item(4, 12).Value = 0.6666666666666666
$ws.Cells.Item(4, 13).Value = 0.1152503333333333
$ws.Cells.Item(4, 14).Value = 0.345751
$ws.Cells.Item(4, 15).Value = 0.03117234266126738
$ws.Cells.Item(4, 16).Value = 0.03117234266126739
$ws.Cells.Item(4, 17).Value = 0.02129580292622223
$ws.Cells.Item(4, 18).Value = 0.191662226336
$ws.Cells.Item(4, 19).Value = 0.03117234266126738
$ws.Cells.Item(4, 20).Value = 0.03117234266126739

# ---- Row 5 (new): Target cluster "sCs" ----
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Vip"
$ws.Cells.Item(5, 3).Value = "Vipr2"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.1847786666666667
$ws.Cells.Item(5, 8).Value = 0.5543360000000001
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.890024666666667
$ws.Cells.Item(5, 14).Value = 8.670074
$ws.Cells.Item(5, 15).Value = 0.7816796412057959
$ws.Cells.Item(5, 16).Value = 0.7816796412057959
$ws.Cells.Item(5, 17).Value = 0.5340149045404444
$ws.Cells.Item(5, 18).Value = 4.806134140864001
$ws.Cells.Item(5, 19).Value = 0.7816796412057959
$ws.Cells.Item(5, 20).Value = 0.7816796412057959
